$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item(1)   # "Sheet1"
$ws2 = $wb.Worksheets.Item(2)   # "Баллон"

# ---------------------------------------------------------------------------
# 1. Updated calculations for selected fabrics: change the fabric thickness
#    and scale inputs on the "Баллон" sheet. Everything downstream
#    (E12, E2:E10, E14, E17, C19, E19, Sheet1!B3/D3/D15/D16/D20, Баллон!E25 …)
#    recalculates automatically.
# ---------------------------------------------------------------------------
$ws2.Range("F12").Value = 1.1000000000000001
$ws2.Range("B17").Value = 0.0001

# ---------------------------------------------------------------------------
# 2. Feather debugging info block (I21/J21) - "красненькая ткань"
# ---------------------------------------------------------------------------
$ws2.Range("I21").Value = "красненькая ткань"
$ws2.Range("J21").Value = 100

# ---------------------------------------------------------------------------
# 3. New "Размеры" input block + "Масштаб"/"Итоговый размер" computations
#    (entered in this order so new shared strings land at the same indices
#    as the authored workbook: Размеры, Итоговый размер, Выкройка, ...)
# ---------------------------------------------------------------------------
$ws2.Range("A28").Value = "Размеры"
$ws2.Range("B28").Value = 1407.21
$ws2.Range("C28").Value = 174.93
$ws2.Range("D28").Value = 79.84
$ws2.Range("E28").Value = 82.38

$ws2.Range("A30").Value = "Итоговый размер"
$ws2.Range("B30").Formula = "=B28*B29"
$ws2.Range("C30").Formula = "=C28*B29"
$ws2.Range("D30").Formula = "=D28*B29"
$ws2.Range("E30").Formula = "=E28*B29"

$ws2.Range("A26").Value = "Выкройка"

$ws2.Range("A32").Value = "Площадь ткани"
$ws2.Range("B32").Formula = "=B30"
$ws2.Range("C32").Formula = "=12*C30"

$ws2.Range("D32").Value = "м"
$ws2.Range("G32").Value = "м2"

$ws2.Range("A29").Value = "Масштаб"
$ws2.Range("B29").Formula = "=(E2*1000)/(B28+D28+E28)"

$ws2.Range("F32").Formula = "=B32/1000 * C32/1000"

# ---------------------------------------------------------------------------
# 4. Refresh the conditional-formatting rules that flag a positive
#    lift margin (Sheet1!D20 and Баллон!E25) with the green fill.
# ---------------------------------------------------------------------------
$cf2 = $ws2.Range("E25").FormatConditions.Item(1)
$cf2.Modify(1, 5, "0")
$cf2.Interior.Color = 5287936

$cf1 = $ws1.Range("D20").FormatConditions.Item(1)
$cf1.Modify(1, 5, "0")
$cf1.Interior.Color = 5287936

# ---------------------------------------------------------------------------
# 5. Restore the last-used selection on each sheet (Баллон stays the active
#    tab, matching the saved view).
# ---------------------------------------------------------------------------
$ws1.Activate()
$ws1.Range("A39").Select()

$ws2.Activate()
$ws2.Range("F13").Select()

Write-Host "done"
